# Commit: "cohort for pos controls added"
#
# The workbook's single sheet ("Metagenome.environmental.1.0") has a
# "Cohort" column (column V). For the positive/negative control samples in
# rows 881-924 the Cohort value was "nan" (blank/not-applicable) and is now
# filled in with the appropriate cohort label, matching the sample id in
# column A for each row:
#   - rows 881-889  (Au14/MC#)                -> MockCommunity
#   - rows 890-905  (18Ja24/ColiGuard|Protexin)-> ColiGuard / Protexin
#   - rows 906-924  (NA/neg.control_#)         -> NegativeControl

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("V881:V889").Value  = "MockCommunity"
$ws.Range("V890").Value       = "ColiGuard"
$ws.Range("V891:V892").Value  = "Protexin"
$ws.Range("V893").Value       = "ColiGuard"
$ws.Range("V894").Value       = "Protexin"
$ws.Range("V895:V896").Value  = "ColiGuard"
$ws.Range("V897:V898").Value  = "Protexin"
$ws.Range("V899").Value       = "ColiGuard"
$ws.Range("V900").Value       = "Protexin"
$ws.Range("V901").Value       = "ColiGuard"
$ws.Range("V902").Value       = "Protexin"
$ws.Range("V903").Value       = "ColiGuard"
$ws.Range("V904").Value       = "Protexin"
$ws.Range("V905").Value       = "ColiGuard"
$ws.Range("V906:V924").Value  = "NegativeControl"

# Leave the view scrolled/selected where the edit ended up, same as the
# author's saved sheet view (selection moves to Y928 after editing the
# Cohort column down through row 924).
[void]$ws.Range("Y928").Select()
